$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries reorder: Austria now ranks (by total cases) ahead of Japon ---
# Row 37 used to be Japon (old data), row 38 used to be Austria (old data).
# After the refresh, Austria (updated figures) takes row 37 and Japon
# (its previous, now-lower-ranked figures) takes row 38.
$ws.Range("A37").Value = "Austria"
$ws.Range("B37").Value = 15997
$ws.Range("C37").Value = 36
$ws.Range("D37").Value = 14304
$ws.Range("E37").Value = 1069
$ws.Range("F37").Value = 55
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 624

$ws.Range("A38").Value = "Japon"
$ws.Range("B38").Value = 15968
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 8531
$ws.Range("E38").Value = 6780
$ws.Range("F38").Value = 249
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 657

# --- Countries reorder: Nepal now ranks ahead of Haiti ---
# Row 142 used to be Haiti (old data), row 143 used to be Nepal (old data).
# After the refresh, Nepal (updated figures) takes row 142 and Haiti
# (its previous figures) takes row 143.
$ws.Range("A142").Value = "Nepal"
$ws.Range("B142").Value = 219
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 33
$ws.Range("E142").Value = 186
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

$ws.Range("A143").Value = "Haiti"
$ws.Range("B143").Value = 219
$ws.Range("C143").Value = 10
$ws.Range("D143").Value = 17
$ws.Range("E143").Value = 184
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 18

# --- Updated timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 10:35"

# --- Updated per-country statistics (no reordering involved) ---
# India
$ws.Range("B15").Value = 74925
$ws.Range("C15").Value = 633
$ws.Range("D15").Value = 24887
$ws.Range("E15").Value = 47602
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 21
$ws.Range("H15").Value = 2436

# Polonia
$ws.Range("B33").Value = 17062
$ws.Range("C33").Value = 141
$ws.Range("D33").Value = 6410
$ws.Range("E33").Value = 9805
$ws.Range("F33").Value = 160
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 847

# Filipinas
$ws.Range("B42").Value = 11618
$ws.Range("C42").Value = 268
$ws.Range("D42").Value = 2251
$ws.Range("E42").Value = 8595
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 21
$ws.Range("H42").Value = 772

# Estonia
$ws.Range("B85").Value = 1751
$ws.Range("C85").Value = 5
$ws.Range("D85").Value = 777
$ws.Range("E85").Value = 913
$ws.Range("F85").Value = 5
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 61

# Sri Lanka
$ws.Range("B104").Value = 889
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 382
$ws.Range("E104").Value = 498
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 9

# Laos
$ws.Range("B190").Value = 19
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 14
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0
